$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Priority" header in D1 (keep its bold style, just remove the text)
$ws.Range("D1").ClearContents()

# Add Complexity numeric values for the first few rows
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 6

# Move the active selection to C4 (also resets the view's top-left scroll position)
$ws.Range("C4").Select()
